$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows for longitude/latitude coordinates
$ws.Range("A10").Value = "LOC_LONGDITUDE"
$ws.Range("B10").Value = -2.8977909999999998

$ws.Range("A11").Value = "LOC_LATITUDE"
$ws.Range("B11").Value = 43.257928

# Update BASE_DIR value (B3) from old repo folder name to new one
$ws.Range("B3").Value = "spacer-hb-framework"

# Update the active cell selection
$ws.Range("F7").Select()
